$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '28.947.30'
$ws.Range("D2").Style = $ws.Range("A1").Style
$ws.Range("E2").Value = '  -1.56%  '
$ws.Range("D3").Value = "'" + '1.910.53'
$ws.Range("D3").Style = $ws.Range("A1").Style
$ws.Range("E3").Value = '  -1.87%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = "'" + '324.97'
$ws.Range("D5").Style = $ws.Range("A1").Style
$ws.Range("E5").Value = '  -0.24%  '
$ws.Range("E6").Value = '  +0.13%  '
$ws.Range("E7").Value = '  -0.91%  '
$ws.Range("D8").Value = "'" + '0.3824'
$ws.Range("D8").Style = $ws.Range("A1").Style
$ws.Range("E8").Value = '  -1.12%  '
$ws.Range("D9").Value = "'" + '0.07729'
$ws.Range("D9").Style = $ws.Range("A1").Style
$ws.Range("E9").Value = '  -1.25%  '
$ws.Range("D10").Value = "'" + '0.9807'
$ws.Range("D10").Style = $ws.Range("A1").Style
$ws.Range("E10").Value = '  +0.30%  '
$ws.Range("E11").Value = '  -2.66%  '
$ws.Range("D12").Value = "'" + '1.904.10'
$ws.Range("D12").Style = $ws.Range("A1").Style
$ws.Range("E12").Value = '  -2.69%  '
$ws.Range("D13").Value = "'" + '6.945'
$ws.Range("D13").Style = $ws.Range("A1").Style
$ws.Range("E13").Value = '  -1.89%  '
$ws.Range("D14").Value = "'" + '5.668'
$ws.Range("D14").Style = $ws.Range("A1").Style
$ws.Range("E14").Value = '  -1.54%  '
$ws.Range("D15").Value = "'" + '0.07040'
$ws.Range("D15").Style = $ws.Range("A1").Style
$ws.Range("E15").Value = '  -0.07%  '
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("E17").Value = '  -3.32%  '
$ws.Range("D18").Value = "'" + '0.000009465'
$ws.Range("D18").Style = $ws.Range("A1").Style
$ws.Range("D19").Value = "'" + '16.69'
$ws.Range("D19").Style = $ws.Range("A1").Style
$ws.Range("E19").Value = '  -2.15%  '
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").Value = "'" + '28.937.56'
$ws.Range("D21").Style = $ws.Range("A1").Style
$ws.Range("E21").Value = '  -1.64%  '
$ws.Range("D22").Value = "'" + '5.327'
$ws.Range("D22").Style = $ws.Range("A1").Style
$ws.Range("E22").Value = '  -2.65%  '
$ws.Range("E23").Value = '  -1.66%  '
$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = "'" + '2.095'
$ws.Range("D24").Style = $ws.Range("A1").Style
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = "'" + '158.93'
$ws.Range("D25").Style = $ws.Range("A1").Style
$ws.Range("E25").Value = '  +1.09%  '
$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").Value = "'" + '19.05'
$ws.Range("D26").Style = $ws.Range("A1").Style
$ws.Range("E26").Value = '  -1.70%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").Value = "'" + '5.670'
$ws.Range("D27").Style = $ws.Range("A1").Style
$ws.Range("E27").Value = '  -1.92%  '
$ws.Range("B28").Value = 'BitcoinCash'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D28").Value = "'" + '117.52'
$ws.Range("D28").Style = $ws.Range("A1").Style
$ws.Range("E28").Value = '  -0.85%  '
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").Value = "'" + '1.855'
$ws.Range("D29").Style = $ws.Range("A1").Style
$ws.Range("E29").Value = '  -0.74%  '
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = "'" + '0.09292'
$ws.Range("D30").Style = $ws.Range("A1").Style
$ws.Range("E30").Value = '  -0.77%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").Value = "'" + '0.8671'
$ws.Range("D31").Style = $ws.Range("A1").Style
$ws.Range("E31").Value = '  +0.25%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = "'" + '5.080'
$ws.Range("D32").Style = $ws.Range("A1").Style
$ws.Range("E32").Value = '  -2.05%  '
$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D33").Value = "'" + '1.252'
$ws.Range("D33").Style = $ws.Range("A1").Style
$ws.Range("E33").Value = '  -4.06%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = "'" + '3.105'
$ws.Range("D34").Style = $ws.Range("A1").Style
$ws.Range("E34").Value = '  -0.56%  '
$ws.Range("B35").Value = 'TrustWalletToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D35").Value = "'" + '1.170'
$ws.Range("D35").Style = $ws.Range("A1").Style
$ws.Range("E35").Value = '  +1.44%  '
$ws.Range("D36").Value = "'" + '0.05726'
$ws.Range("D36").Style = $ws.Range("A1").Style
$ws.Range("E36").Value = '  -0.73%  '
$ws.Range("B37").Value = 'Frax'
$ws.Range("C37").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").Value = "'" + '1.002'
$ws.Range("D37").Style = $ws.Range("A1").Style
$ws.Range("E37").Value = '  +0.09%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = "'" + '0.02042'
$ws.Range("D38").Style = $ws.Range("A1").Style
$ws.Range("E38").Value = '  -2.03%  '
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").Value = "'" + '0.5496'
$ws.Range("D39").Style = $ws.Range("A1").Style
$ws.Range("E39").Value = '  -2.97%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").Value = "'" + '7.409'
$ws.Range("D40").Style = $ws.Range("A1").Style
$ws.Range("E40").Value = '  -3.79%  '
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = "'" + '0.1754'
$ws.Range("D41").Style = $ws.Range("A1").Style
$ws.Range("E41").Value = '  -1.68%  '
$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D42").Value = "'" + '2.857'
$ws.Range("D42").Style = $ws.Range("A1").Style
$ws.Range("E42").Value = '  +4.92%  '
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = "'" + '9.317'
$ws.Range("D43").Style = $ws.Range("A1").Style
$ws.Range("E43").Value = '  -1.26%  '
$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").Value = "'" + '0.5186'
$ws.Range("D44").Style = $ws.Range("A1").Style
$ws.Range("E44").Value = '  -2.03%  '
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").Value = "'" + '0.06905'
$ws.Range("D45").Style = $ws.Range("A1").Style
$ws.Range("E45").Value = '  +0.46%  '
$ws.Range("E46").Value = '  -3.84%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = "'" + '2.079'
$ws.Range("D47").Style = $ws.Range("A1").Style
$ws.Range("E47").Value = '  -0.80%  '
$ws.Range("B48").Value = 'PEPE'
$ws.Range("C48").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D48").Value = "'" + '0.000002581'
$ws.Range("D48").Style = $ws.Range("A1").Style
$ws.Range("E48").Value = '  -11.02%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = "'" + '1.781'
$ws.Range("D49").Style = $ws.Range("A1").Style
$ws.Range("E49").Value = '  -2.05%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = "'" + '110.57'
$ws.Range("D50").Style = $ws.Range("A1").Style
$ws.Range("E50").Value = '  -0.86%  '
$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").Value = "'" + '0.2875'
$ws.Range("D51").Style = $ws.Range("A1").Style
$ws.Range("E51").Value = '  -4.24%  '